# Append new job listing row (2026-01-24 01:57:29 JST crawl) into the
# "ランサーズ" sheet, refresh the crawl timestamp on all existing rows,
# and keep the F-column hyperlinks pointing at the right cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-24 01:57:29"

# Insert a fresh row 4 (existing rows 4-15 shift down to 5-16).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new listing.
$ws.Cells.Item(4, 1).Value = $newTimestamp
$ws.Cells.Item(4, 2).Value = "【急募】Amazon SP-API 自動化開発者を探しています"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5477903"
$ws.Cells.Item(4, 7).Value = 328
$ws.Cells.Item(4, 8).Value = "🔥API ◆開発,自動化"

# Refresh the "fetched at" timestamp on every other data row (2-3, 5-16).
$lastRow = 16
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Rebuild the F-column hyperlinks so ref cells line up with the shifted
# rows (row insert does not relocate existing hyperlinks/relationships).
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le $lastRow; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
